$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.124.25'
$ws.Range('E2').Value = '  +1.56%  '

$ws.Range('D3').Value = '3.218.91'
$ws.Range('E3').Value = '  +1.40%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.89'
$ws.Range('E5').Value = '  +4.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.97'
$ws.Range('E6').Value = '  +1.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').Value = '3.218.67'
$ws.Range('E8').Value = '  +1.63%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +1.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  -1.06%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.15'
$ws.Range('E11').Value = '  -0.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.512'
$ws.Range('E12').Value = '  +1.13%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('E13').Value = '  -1.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.58'
$ws.Range('E14').Value = '  +1.37%  '

$ws.Range('D15').Value = '3.734.52'
$ws.Range('E15').Value = '  +1.09%  '

$ws.Range('D16').Value = '66.088.59'
$ws.Range('E16').Value = '  +1.35%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.41'
$ws.Range('E17').Value = '  +3.08%  '

$ws.Range('D18').Value = '3.206.14'
$ws.Range('E18').Value = '  +0.58%  '

$ws.Range('E19').Value = '  +0.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '513.57'
$ws.Range('E20').Value = '  -0.02%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.83'
$ws.Range('E21').Value = '  +6.16%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.737'
$ws.Range('E22').Value = '  +0.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.19'
$ws.Range('E23').Value = '  -1.31%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.00'
$ws.Range('E24').Value = '  +2.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.48'
$ws.Range('E25').Value = '  +0.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.34'
$ws.Range('E27').Value = '  +2.78%  '

$ws.Range('E28').Value = '  +3.65%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +2.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.88'
$ws.Range('E30').Value = '  +3.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.82'
$ws.Range('E31').Value = '  +7.82%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.24'
$ws.Range('E32').Value = '  +0.42%  '

$ws.Range('E33').Value = '  +0.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.13%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.65'
$ws.Range('E35').Value = '  -0.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.58'
$ws.Range('E36').Value = '  -0.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0926'
$ws.Range('E37').Value = '  +3.00%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '487.96'
$ws.Range('E38').Value = '  +1.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0422'
$ws.Range('E39').Value = '  +0.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.01'
$ws.Range('E40').Value = '  -2.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.88'
$ws.Range('E41').Value = '  +2.71%  '

$ws.Range('D42').Value = '3.031.71'
$ws.Range('E42').Value = '  -2.37%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.120'
$ws.Range('E43').Value = '  +0.01%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.295'
$ws.Range('E44').Value = '  +2.35%  '

$ws.Range('D45').Value = '0.0₃0649'
$ws.Range('E45').Value = '  +8.46%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.46'
$ws.Range('E46').Value = '  +0.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '29.18'
$ws.Range('E47').Value = '  -0.85%  '

$ws.Range('E48').Value = '  +0.07%  '

$ws.Range('E49').Value = '  +0.66%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.33'
$ws.Range('E50').Value = '  +1.71%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.79'
$ws.Range('E51').Value = '  -1.51%  '
